$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.127400756369923
$ws.Range("C2").Value = 0.2477383497348171
$ws.Range("D2").Value = 0.07879849045720277
$ws.Range("E2").Value = 0.08368605951984165
$ws.Range("G2").Value = 1.101277027665645
$ws.Range("H2").Value = 0.9832047830831527
$ws.Range("M2").Value = 0.4077578897561907
$ws.Range("N2").Value = 1.27710018740224
# Row 3
$ws.Range("B3").Value = 1.008614740519931
$ws.Range("C3").Value = 0.2157094906861801
$ws.Range("D3").Value = 0.07141371160065546
$ws.Range("E3").Value = 0.07895873490489436
$ws.Range("G3").Value = 1.05495773986928
$ws.Range("H3").Value = 0.9679947586264177
$ws.Range("M3").Value = 0.36724502438922
$ws.Range("N3").Value = 1.292190060109959
# Row 4
$ws.Range("B4").Value = 0.9361049630625757
$ws.Range("C4").Value = 0.196050933656096
$ws.Range("D4").Value = 0.06692133839526093
$ws.Range("E4").Value = 0.07611931704547814
$ws.Range("G4").Value = 1.027386304530808
$ws.Range("H4").Value = 0.959310425863606
$ws.Range("M4").Value = 0.3425646113710528
$ws.Range("N4").Value = 1.30201478526979
# Row 5
$ws.Range("B5").Value = 0.9066618503301811
$ws.Range("C5").Value = 0.1880411862231028
$ws.Range("D5").Value = 0.06510102353043123
$ws.Range("E5").Value = 0.07497787864892658
$ws.Range("G5").Value = 1.016366830089396
$ws.Range("H5").Value = 0.9559351628988964
$ws.Range("M5").Value = 0.3325551625862175
$ws.Range("N5").Value = 1.30615855455666
# Row 6
$ws.Range("B6").Value = 0.9017791482211237
$ws.Range("C6").Value = 0.186711234555105
$ws.Range("D6").Value = 0.06479938249943018
$ws.Range("E6").Value = 0.07478928178827715
$ws.Range("G6").Value = 1.014550038724138
$ws.Range("H6").Value = 0.9553845608666052
$ws.Range("M6").Value = 0.330895974706884
$ws.Range("N6").Value = 1.306855070185428
# Row 7
$ws.Range("B7").Value = 0.9357074592289223
$ws.Range("C7").Value = 0.1959429069676162
$ws.Range("D7").Value = 0.06689674725814143
$ws.Range("E7").Value = 0.07610386016581927
$ws.Range("G7").Value = 1.027236820137034
$ws.Range("H7").Value = 0.959264244451731
$ws.Range("M7").Value = 0.342429427224836
$ws.Range("N7").Value = 1.302070103049864
# Row 8
$ws.Range("B8").Value = 1.08635397168274
$ws.Range("C8").Value = 0.2366929374947233
$ws.Range("D8").Value = 0.07624341190407335
$ws.Range("E8").Value = 0.08204281906232325
$ws.Range("G8").Value = 1.085124379230706
$ws.Range("H8").Value = 0.9778239037185017
$ws.Range("M8").Value = 0.3937480808051319
$ws.Range("N8").Value = 1.282186710727522
# Row 9
$ws.Range("B9").Value = 1.385240860850899
$ws.Range("C9").Value = 0.3166908013192824
$ws.Range("D9").Value = 0.09491364527143276
$ws.Range("E9").Value = 0.09420190060018996
$ws.Range("G9").Value = 1.205648888094373
$ws.Range("H9").Value = 1.019460377253495
$ws.Range("M9").Value = 0.4959756231912991
$ws.Range("N9").Value = 1.247659971920307
# Row 10
$ws.Range("B10").Value = 1.607103182306616
$ws.Range("C10").Value = 0.3755648734850183
$ws.Range("D10").Value = 0.1088526151265228
$ws.Range("E10").Value = 0.1034647598186993
$ws.Range("G10").Value = 1.29863812211272
$ws.Range("H10").Value = 1.053315043973925
$ws.Range("M10").Value = 0.5721262307126693
$ws.Range("N10").Value = 1.225049145398756
# Row 11
$ws.Range("B11").Value = 1.708562395995273
$ws.Range("C11").Value = 0.4023803639837524
$ws.Range("D11").Value = 0.1152451002688366
$ws.Range("E11").Value = 0.107754003680931
$ws.Range("G11").Value = 1.341942667360655
$ws.Range("H11").Value = 1.069440679251102
$ws.Range("M11").Value = 0.6070123797577196
$ws.Range("N11").Value = 1.215368818418078
# Row 12
$ws.Range("B12").Value = 1.747061430014298
$ws.Range("C12").Value = 0.4125402372526992
$ws.Range("D12").Value = 0.1176734091960014
$ws.Range("E12").Value = 0.1093893796649112
$ws.Range("G12").Value = 1.3584880999787
$ws.Range("H12").Value = 1.075652475245477
$ws.Range("M12").Value = 0.6202592693230571
$ws.Range("N12").Value = 1.211790851892701
# Row 13
$ws.Range("B13").Value = 1.738766444771954
$ws.Range("C13").Value = 0.4103518751627462
$ws.Range("D13").Value = 0.1171500876440206
$ws.Range("E13").Value = 0.1090366726908485
$ws.Range("G13").Value = 1.354918166888183
$ws.Range("H13").Value = 1.074309949598074
$ws.Range("M13").Value = 0.6174046849277914
$ws.Range("N13").Value = 1.212557516448832
# Row 14
$ws.Range("B14").Value = 1.711728147202678
$ws.Range("C14").Value = 0.4032161104207717
$ws.Range("D14").Value = 0.1154447249814581
$ws.Range("E14").Value = 0.1078883223630029
$ws.Range("G14").Value = 1.343300909720512
$ws.Range("H14").Value = 1.069949609259595
$ws.Range("M14").Value = 0.6081014764973389
$ws.Range("N14").Value = 1.215072694008519
# Row 15
$ws.Range("B15").Value = 1.695176709412976
$ws.Range("C15").Value = 0.3988459743659405
$ws.Range("D15").Value = 0.1144011385060537
$ws.Range("E15").Value = 0.107186382549024
$ws.Range("G15").Value = 1.336204220188876
$ws.Range("H15").Value = 1.067292530345185
$ws.Range("M15").Value = 0.6024077466360893
$ws.Range("N15").Value = 1.216624763194673
# Row 16
$ws.Range("B16").Value = 1.600483477730904
$ws.Range("C16").Value = 0.3738131433019589
$ws.Range("D16").Value = 0.10843590581608
$ws.Range("E16").Value = 0.1031859905482762
$ws.Range("G16").Value = 1.295828484314683
$ws.Range("H16").Value = 1.052275879287492
$ws.Range("M16").Value = 0.569851355475123
$ws.Range("N16").Value = 1.225694020373254
# Row 17
$ws.Range("B17").Value = 1.542530150814059
$ws.Range("C17").Value = 0.3584653054201112
$ws.Range("D17").Value = 0.1047897781482163
$ws.Range("E17").Value = 0.100751415374539
$ws.Range("G17").Value = 1.271318131483497
$ws.Range("H17").Value = 1.043250123867324
$ws.Range("M17").Value = 0.5499425203147439
$ws.Range("N17").Value = 1.231413281411413
# Row 18
$ws.Range("B18").Value = 1.509246969901767
$ws.Range("C18").Value = 0.3496407385229077
$ws.Range("D18").Value = 0.1026974647627839
$ws.Range("E18").Value = 0.09935820416716723
$ws.Range("G18").Value = 1.257314696223716
$ws.Range("H18").Value = 1.038126918608214
$ws.Range("M18").Value = 0.5385144678701437
$ws.Range("N18").Value = 1.234759794591248
# Row 19
$ws.Range("B19").Value = 1.497986390997426
$ws.Range("C19").Value = 0.3466534058648563
$ws.Range("D19").Value = 0.1019898692134262
$ws.Range("E19").Value = 0.09888769662097019
$ws.Range("G19").Value = 1.252589484002073
$ws.Range("H19").Value = 1.036403963918559
$ws.Range("M19").Value = 0.5346490412619573
$ws.Range("N19").Value = 1.23590262571873
# Row 20
$ws.Range("B20").Value = 1.548694188128025
$ws.Range("C20").Value = 0.3600987832423357
$ws.Range("D20").Value = 0.1051774120348199
$ws.Range("E20").Value = 0.1010098440982787
$ws.Range("G20").Value = 1.273917524408688
$ws.Range("H20").Value = 1.044203867513005
$ws.Range("M20").Value = 0.5520594635905809
$ws.Range("N20").Value = 1.230798557631026
# Row 21
$ws.Range("B21").Value = 1.71966780044329
$ws.Range("C21").Value = 0.4053119050854548
$ws.Range("D21").Value = 0.115945423122767
$ws.Range("E21").Value = 0.1082253163702021
$ws.Range("G21").Value = 1.346709170361805
$ws.Range("H21").Value = 1.071227478272334
$ws.Range("M21").Value = 0.6108330617100535
$ws.Range("N21").Value = 1.214331537891717
# Row 22
$ws.Range("B22").Value = 1.831868223705044
$ws.Range("C22").Value = 0.4348932387536593
$ws.Range("D22").Value = 0.1230274117710763
$ws.Range("E22").Value = 0.113006062081638
$ws.Range("G22").Value = 1.395140468074857
$ws.Range("H22").Value = 1.089503665218814
$ws.Range("M22").Value = 0.6494568699974366
$ws.Range("N22").Value = 1.204081271807084
# Row 23
$ws.Range("B23").Value = 1.771942050512621
$ws.Range("C23").Value = 0.4191019903253164
$ws.Range("D23").Value = 0.1192434862966678
$ws.Range("E23").Value = 0.1104484480140968
$ws.Range("G23").Value = 1.369212397578025
$ws.Range("H23").Value = 1.079692702155256
$ws.Range("M23").Value = 0.628822878630757
$ws.Range("N23").Value = 1.209504960421889
# Row 24
$ws.Range("B24").Value = 1.545907318120214
$ws.Range("C24").Value = 0.3593602907396303
$ws.Range("D24").Value = 0.1050021506408143
$ws.Range("E24").Value = 0.1008929883760388
$ws.Range("G24").Value = 1.272742065301514
$ws.Range("H24").Value = 1.043772475075173
$ws.Range("M24").Value = 0.5511023381652649
$ws.Range("N24").Value = 1.23107629229667
# Row 25
$ws.Range("B25").Value = 1.303996726228547
$ws.Range("C25").Value = 0.2950351608429571
$ws.Range("D25").Value = 0.08982488186347837
$ws.Range("E25").Value = 0.09085590723140058
$ws.Range("G25").Value = 1.172276381448285
$ws.Range("H25").Value = 1.007628414890377
$ws.Range("M25").Value = 0.4681424006549548
$ws.Range("N25").Value = 1.256518903598845
